# Temple Plan V1 - info/scene excel file follow-up edits
# Slide 4 (index 4): Waiting Room label, Crafting Room label + line break,
#                     and 3 connector line colors (accent6 -> accent2)
# Slide 5 (index 5): Prayer Room label, Room 1/2/3 relabeled with old-room notes

$p = $ppt.ActivePresentation

# ----- Slide 4 -----
$s4 = $p.Slides.Item(4)

# "Visitor/Waiting Room" -> "Waiting Room"
$s4.Shapes.Item("Rounded Rectangle 1").TextFrame.TextRange.Text = "Waiting Room"

# "Crafting (art) Room" -> "Crafting  Room" + a manual (soft) line break, leaving
# a trailing empty line in the box
$s4.Shapes.Item("Rounded Rectangle 73").TextFrame.TextRange.Text = "Crafting  Room" + [char]11

# Connector line colors: accent6 -> accent2
$s4.Shapes.Item("Curved Connector 75").Line.ForeColor.ObjectThemeColor = 6
$s4.Shapes.Item("Curved Connector 87").Line.ForeColor.ObjectThemeColor = 6
$s4.Shapes.Item("Curved Connector 90").Line.ForeColor.ObjectThemeColor = 6

# ----- Slide 5 -----
$s5 = $p.Slides.Item(5)

# "Common Prayer Room" -> "Prayer Room"
$s5.Shapes.Item("Rounded Rectangle 13").TextFrame.TextRange.Text = "Prayer Room"

# "Room 1" -> "RO's Room" + line break + "(EO's old Room)"
$s5.Shapes.Item("Rounded Rectangle 50").TextFrame.TextRange.Text = "RO" + [char]8217 + "s Room" + [char]11 + "(EO" + [char]8217 + "s old Room)"

# "Room 2" -> "SO's Room" + line break + "(RO's old Room)"
$s5.Shapes.Item("Rounded Rectangle 51").TextFrame.TextRange.Text = "SO" + [char]8217 + "s Room" + [char]11 + "(RO" + [char]8217 + "s old Room)"

# "Room 3" -> "LO's Room"
$s5.Shapes.Item("Rounded Rectangle 52").TextFrame.TextRange.Text = "LO" + [char]8217 + "s Room"
